$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 569.062204
$ws.Range("D2").Value = 54.631182

$ws.Range("B3").Value = 144.664285
$ws.Range("D3").Value = 6.94404
$ws.Range("E3").Value = 0.001187

$ws.Range("B4").Value = 2322.86518
$ws.Range("C4").Value = 223

$ws.Range("G5").Value = -1.286
$ws.Range("H5").Value = -2.695544
$ws.Range("I5").Value = 0.123544
$ws.Range("J5").Value = 0.081815

$ws.Range("G6").Value = 0.444548
$ws.Range("H6").Value = -1.073233
$ws.Range("I6").Value = 1.962329
$ws.Range("J6").Value = 0.76895

$ws.Range("G7").Value = 1.730548
$ws.Range("H7").Value = 0.58309
$ws.Range("I7").Value = 2.878006
$ws.Range("J7").Value = 0.00132
